$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D (Price) to Text format so values like "21.50" or
# "1.00" are written verbatim instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "59.230.85"
$ws.Range("E2").Value = "  -2.40%  "
$ws.Range("D3").Value = "2.636.56"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "515.86"
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("D6").Value = "149.53"
$ws.Range("E6").Value = "  -3.01%  "
$ws.Range("D7").Value = "0.995"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").Value = "0.579"
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("D9").Value = "2.665.68"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").Value = "6.52"
$ws.Range("E10").Value = "  +0.68%  "
$ws.Range("D11").Value = "0.107"
$ws.Range("E11").Value = "  -1.65%  "
$ws.Range("D12").Value = "0.342"
$ws.Range("E12").Value = "  -2.06%  "
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("D14").Value = "3.098.06"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15").Value = "59.062.10"
$ws.Range("E15").Value = "  -2.68%  "
$ws.Range("D16").Value = "21.50"
$ws.Range("E16").Value = "  -1.78%  "
$ws.Range("D17").Value = "0.0000140"
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").Value = "2.648.59"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").Value = "4.64"
$ws.Range("E19").Value = "  -2.29%  "
$ws.Range("D20").Value = "348.35"
$ws.Range("E20").Value = "  -1.31%  "
$ws.Range("D21").Value = "10.62"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").Value = "6.29"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "61.20"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").Value = "0.427"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "2.735.55"
$ws.Range("E26").Value = "  -1.34%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "0.162"
$ws.Range("E27").Value = "  -2.52%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "0.994"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("D29").Value = "0.0₃0839"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("D30").Value = "7.17"
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("D32").Value = "6.44"
$ws.Range("E32").Value = "  +4.54%  "
$ws.Range("D33").Value = "19.25"
$ws.Range("E33").Value = "  -0.63%  "
$ws.Range("D34").Value = "1.58"
$ws.Range("E34").Value = "  -2.43%  "
$ws.Range("D35").Value = "149.20"
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("D36").Value = "1.05"
$ws.Range("E36").Value = "  +18.54%  "
$ws.Range("D37").Value = "4.10"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "1.16"
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("D39").Value = "0.891"
$ws.Range("E39").Value = "  -1.79%  "
$ws.Range("D40").Value = "36.37"
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("D41").Value = "1.45"
$ws.Range("E41").Value = "  -0.81%  "
$ws.Range("D42").Value = "3.71"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").Value = "295.79"
$ws.Range("E43").Value = "  -3.17%  "
$ws.Range("D44").Value = "0.628"
$ws.Range("E44").Value = "  -1.70%  "
$ws.Range("D45").Value = "0.100"
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "19.93"
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").Value = "0.993"
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("D48").Value = "0.0548"
$ws.Range("E48").Value = "  -1.91%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "0.0234"
$ws.Range("E49").Value = "  -2.63%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "4.80"
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "19.21"
$ws.Range("E51").Value = "  +0.32%  "

# Restore the original (unstyled) cell style now that the exact text is stored.
$ws.Range("D2:D51").Style = "Normal"
